# Remove the two MQDSS benchmark rows ("MQDSS-31-48" at row 14 and
# "MQDSS-31-64" at row 20) from the "Data" table. Deleting the
# higher-numbered row first keeps row 14's index valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Rows("20").Delete()
$ws.Rows("14").Delete()

# The ListObject (Algorithm_benchmark_100) range/autofilter/sort state
# shrinks automatically when the rows are removed, but the chart series,
# the hidden ExternalData_1 defined name, and the chart's drawing anchor
# still reference the old A1:L62 extent, so bring those up to date by hand.

$co = $ws.ChartObjects()
$chart = $co.Item(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Data!`$I`$1,Data!`$A`$2:`$A`$60,Data!`$I`$2:`$I`$60,1)"

$nm = $wb.Names.Item("ExternalData_1")
$nm.RefersTo = "=Data!`$A`$1:`$L`$60"

# The chart's graphic frame is anchored with a fixed bottom-right cell
# (row 41 in the original file); shrink it by the 2 deleted rows (30pt
# at the sheet's default 15pt row height) so it again ends at row 39.
$chartObj = $co.Item(1)
$chartObj.Height = $chartObj.Height - 30

# Update the active selection to match the post-edit state.
$ws.Range("A26").Select()
